$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new column before column F. This shifts the old F,G,H,I
#    columns (for all 128 rows) one column to the right, becoming
#    G,H,I,J, and also widens the sheet's used range to column J.
# ------------------------------------------------------------------
$ws.Columns("F:F").Insert()

# ------------------------------------------------------------------
# 2. The table in rows 2-11 (columns C through J) was substantially
#    reorganised by hand (a new "Estimated days" column was added and
#    a row of data was re-shuffled). Rather than rely on the partial
#    shift performed by the column insert above, explicitly write out
#    the exact final contents of every affected cell.
# ------------------------------------------------------------------

# Row 1 headers
$ws.Range("F1").Value = "Estimated days"
$ws.Range("H1").Value = "ESTIMATED REMOVAL (m)"

# Row 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 23
$ws.Range("F2").Value = 23
$ws.Range("G2").Value = "DATE"
$ws.Range("H2").Value = "DOY"
$ws.Range("I2").Value = "level"
$ws.Range("J2").Value = "removal amount (m3)"

# Row 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Formula = '=SLOPE(B7:B35,A7:A35)'
$ws.Range("E3").Value = 66
$ws.Range("F3").Formula = '=A35-F2'
$ws.Range("G3").Value = "Jan 23 - Jan 25"
$ws.Range("H3").Value = "23-25"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("I3").Formula = '=B5-B7'
$ws.Range("J3").Formula = '=I3*D$9'

# Row 4
$ws.Range("C4").Value = 3
$ws.Range("D4").Formula = '=SLOPE(B36:B59,A36:A59)'
$ws.Range("E4").Value = 24
$ws.Range("F4").Formula = '=A59-F2-F3'
$ws.Range("G4").Value = "Mar 31 - Apr 1"
$ws.Range("H4").Value = "90-91"
$ws.Range("I4").Formula = '=B35-B36'
$ws.Range("J4").Formula = '=I4*D$9+D4*D$9*1'

# Row 5
$ws.Range("C5").Value = 4
$ws.Range("D5").Formula = '=SLOPE(B62:B104,A62:A104)'
$ws.Range("E5").Value = 89
$ws.Range("F5").Formula = '=A104-SUM(F2:F4)'
$ws.Range("G5").Value = "Apr 23 - Apr 27"
$ws.Range("H5").Value = "113-117"
$ws.Range("I5").Formula = '=B59-B62'
$ws.Range("J5").Formula = '=I5*D$9+D5*D$9*4'

# Row 6
$ws.Range("C6").Value = 5
$ws.Range("D6").Formula = '=SLOPE(B105:B108,A105:A108)'
$ws.Range("E6").Value = 50
$ws.Range("F6").Formula = '=A108-SUM(F2:F5)'
$ws.Range("G6").Value = "Jul 24 - Aug 29"
$ws.Range("H6").Value = "205-241"
$ws.Range("I6").Formula = '=B104-B105'
$ws.Range("J6").Formula = '=I6*D$9+D6*D$9*36'

# Row 7
$ws.Range("C7").Value = 6
$ws.Range("D7").Formula = '=SLOPE(B111:B128,A111:A128)'
$ws.Range("E7").Value = 60
$ws.Range("F7").Formula = '=A128-SUM(F2:F6)+ (366-A128)'
$ws.Range("G7").Value = "Oct 16 - Oct 21"
$ws.Range("H7").Value = "289-294"
$ws.Range("I7").Formula = '=B108-B111'
$ws.Range("J7").Formula = '=I7*D$9+D7*D$9*5'

# Row 8
$ws.Range("C8").Value = "Avg. (m)"
$ws.Range("D8").Formula = '=SUMPRODUCT(D3:D7,E3:E7)/SUM(E3:E7)'
$ws.Range("E8").ClearContents()
$ws.Range("F8").Formula = '=SUM(F2:F7)'
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()

# Row 9
$ws.Range("C9").Value = "tank area"
$ws.Range("D9").Formula = '=20^2*PI()'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").ClearContents()
$ws.Range("E9").NumberFormat = "General"
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = "total removal"
$ws.Range("I9").ClearContents()
$ws.Range("J9").Formula = '=SUM(J3:J7)'
$ws.Range("J9").NumberFormat = "General"

# Row 10
$ws.Range("C10").Value = "daily manure input (m3)"
$ws.Range("D10").Formula = '=D8*D9'
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()

# Row 11
$ws.Range("C11").Value = "annual manure input (m3)"
$ws.Range("D11").Formula = '=D8*D9*366'
$ws.Range("E11").ClearContents()

# ------------------------------------------------------------------
# 3. Re-apply the number formats that belonged to the old D9/E9 (now
#    D10/E10) and D10/E10 (now D11/E11) pairs so the shifted cells
#    keep their original look (0.00 and 0 respectively).
# ------------------------------------------------------------------
$ws.Range("D10:E10").NumberFormat = "0.00"
$ws.Range("F9").NumberFormat = "0.00"
$ws.Range("D11:E11").NumberFormat = "0"
$ws.Range("F10").NumberFormat = "0"

# ------------------------------------------------------------------
# 4. Shift the chart so it keeps sitting to the right of the table,
#    now that an extra column has been inserted in front of it.
# ------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$colFWidth = $ws.Columns("F").Width
$chartObj.Left = $chartObj.Left + $colFWidth

# ------------------------------------------------------------------
# 5. Restore the selected cell shown when the workbook is reopened.
# ------------------------------------------------------------------
$ws.Range("K18").Select()
